$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) to the maternal-mortality table, mirroring
# column R (the "2021" column) for formatting, and filling in the new
# per-region values for row 4 (year header) through row 14 (last data row).

$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Bring column S's formatting in line with column R for the table rows
# (thin/medium border row 3, header row 4, and the data rows 5-14).
Copy-Format "R3" "S3"
Copy-Format "R4" "S4"
Copy-Format "R5" "S5"
Copy-Format "R6" "S6"
Copy-Format "R7" "S7"
Copy-Format "R8" "S8"
Copy-Format "R9" "S9"
Copy-Format "R10" "S10"
Copy-Format "R11" "S11"
Copy-Format "R12" "S12"
Copy-Format "R13" "S13"
Copy-Format "R14" "S14"

# New "2022" year header.
$ws.Range("S4").Value2 = 2022

# New data values per region for 2022.
$ws.Range("S5").Value2 = 27.292394741221504
$ws.Range("S6").Value2 = 36.613942589338023
$ws.Range("S7").Value2 = 14.18691257315127
$ws.Range("S8").Value2 = 55.377118174770182
$ws.Range("S9").Value2 = 42.247570764681029
$ws.Range("S10").Value2 = 30.18817294468856
$ws.Range("S11").Value2 = 97.03085581214826
$ws.Range("S12").Value2 = 25.2
$ws.Range("S13").Value2 = 21.849963583394029
$ws.Range("S14").Value2 = "-"

# Match the author's leftover selection state after filling the column.
$ws.Range("S16").Select() | Out-Null
